$d = $word.ActiveDocument

# --- Change 1: "Data is transferred via the Internet not through a continuous
# stream, but rather in millions of fundamental units known as Packets:"
# becomes three runs: "Data is transferred via the Internet" + "," +
# " not through a continuous stream, but rather in millions of fundamental
# units known as Packets:"
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Data is transferred via the Internet not through a continuous stream*") {
        $r = $p.Range
        $start = $r.Start
        $firstPart = "Data is transferred via the Internet"
        $len1 = $firstPart.Length

        # Insert the comma right after "Internet" (before the following space)
        $insPos = $start + $len1
        $insRng = $d.Range($insPos, $insPos)
        $insRng.InsertAfter(",")

        # Force a run split between "Internet" and the new comma by touching
        # (and reverting) formatting on the first segment.
        $rng1 = $d.Range($start, $start + $len1)
        $rng1.Font.Bold = 1
        $rng1.Font.Bold = 0

        # Force a run split between the comma and the remaining text.
        $rng2 = $d.Range($start + $len1, $start + $len1 + 1)
        $rng2.Font.Bold = 1
        $rng2.Font.Bold = 0

        break
    }
}

# --- Change 2: fix "it's" -> "its" typo in the Fragmentation bullet
$apos = [char]0x2019
$d.Content.Find.Execute("it" + $apos + "s addressing information", $true, $false, $false, $false, $false, $true, 1, $false, "its addressing information", 2)
